$wb = $excel.ActiveWorkbook

# --- Update "sets" sheet: set_id=3 home_points 12 -> 15 ---
$wsSets = $wb.Worksheets.Item("sets")
$wsSets.Cells.Item(4, 4).Value = 15

# --- Add new rows to "rallies" sheet ---
$wsRallies = $wb.Worksheets.Item("rallies")

# Row 76 (rally_id 75)
$wsRallies.Cells.Item(76, 1).Value = 75
$wsRallies.Cells.Item(76, 2).Value = 1
$wsRallies.Cells.Item(76, 3).Value = 3
$wsRallies.Cells.Item(76, 4).Value = 13
$wsRallies.Cells.Item(76, 5).Value = "NOS"
$wsRallies.Cells.Item(76, 6).Value = ""
$wsRallies.Cells.Item(76, 7).Value = 4
$wsRallies.Cells.Item(76, 8).Value = "LOB"
$wsRallies.Cells.Item(76, 9).Value = "PONTO"
$wsRallies.Cells.Item(76, 10).Value = "NOS"
$wsRallies.Cells.Item(76, 11).Value = 13
$wsRallies.Cells.Item(76, 12).Value = 0
$wsRallies.Cells.Item(76, 13).Value = "1 4 lob"
$wsRallies.Cells.Item(76, 14).Value = "FRENTE"
$wsRallies.Cells.Item(76, 15).Value = "FRENTE"
$wsRallies.Cells.Item(76, 16).Value = "FRENTE"

# Row 77 (rally_id 76)
$wsRallies.Cells.Item(77, 1).Value = 76
$wsRallies.Cells.Item(77, 2).Value = 1
$wsRallies.Cells.Item(77, 3).Value = 3
$wsRallies.Cells.Item(77, 4).Value = 14
$wsRallies.Cells.Item(77, 5).Value = "NOS"
$wsRallies.Cells.Item(77, 6).Value = ""
$wsRallies.Cells.Item(77, 7).Value = 5
$wsRallies.Cells.Item(77, 8).Value = "LOB"
$wsRallies.Cells.Item(77, 9).Value = "PONTO"
$wsRallies.Cells.Item(77, 10).Value = "NOS"
$wsRallies.Cells.Item(77, 11).Value = 14
$wsRallies.Cells.Item(77, 12).Value = 0
$wsRallies.Cells.Item(77, 13).Value = "1 5 lob"
$wsRallies.Cells.Item(77, 14).Value = "FRENTE"
$wsRallies.Cells.Item(77, 15).Value = "FRENTE"
$wsRallies.Cells.Item(77, 16).Value = "FRENTE"

# Row 78 (rally_id 77)
$wsRallies.Cells.Item(78, 1).Value = 77
$wsRallies.Cells.Item(78, 2).Value = 1
$wsRallies.Cells.Item(78, 3).Value = 3
$wsRallies.Cells.Item(78, 4).Value = 15
$wsRallies.Cells.Item(78, 5).Value = "NOS"
$wsRallies.Cells.Item(78, 6).Value = ""
$wsRallies.Cells.Item(78, 7).Value = 5
$wsRallies.Cells.Item(78, 8).Value = "SEGUNDA"
$wsRallies.Cells.Item(78, 9).Value = "PONTO"
$wsRallies.Cells.Item(78, 10).Value = "NOS"
$wsRallies.Cells.Item(78, 11).Value = 15
$wsRallies.Cells.Item(78, 12).Value = 0
$wsRallies.Cells.Item(78, 13).Value = "1 5 seg"
$wsRallies.Cells.Item(78, 14).Value = "FRENTE"
$wsRallies.Cells.Item(78, 15).Value = "FRENTE"
$wsRallies.Cells.Item(78, 16).Value = "FRENTE"
